# Updates the "Avverkningsanmälningar" sheet:
#  - bumps the "Förändrad" date (column C) from 2026-03-01 (46078) to
#    2026-03-02 (46079) for every data row (2-33)
#  - re-orders the record rows 5-33 into a new sequence (the whole row's
#    contents - designation, dates, area, species counts, species names,
#    and the six HYPERLINK formulas - travel together as a unit)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Snapshot the current contents of rows 5-33 (29 rows) before any
#    writes happen, so later writes in this script don't clobber data
#    we still need to read.
# ---------------------------------------------------------------------
$oldVals = $ws.Range("A5:R33").Value2          # columns A..R (18 cols)
$oldFormulas = $ws.Range("S5:Y33").Formula     # columns S..Y (7 cols)

# Permutation: for new row (5..33), which old row (5..33) supplies the data.
# (row 7 keeps its own data - it is its own source)
$srcRow = @(6,5,7,25,30,31,28,11,14,9,13,17,21,29,26,8,19,24,10,33,12,18,20,15,16,27,32,23,22)

$rowCount = $srcRow.Length   # 29 (rows 5..33)
$colCountAR = 18             # A..R
$colCountSY = 7              # S..Y

$newVals = New-Object 'object[,]' $rowCount,$colCountAR
$newFormulas = New-Object 'object[,]' $rowCount,$colCountSY

for ($i = 0; $i -lt $rowCount; $i++) {
    $srcIdx = $srcRow[$i] - 5   # 0-based index into the snapshot arrays (row 5 -> 0)

    for ($c = 1; $c -le $colCountAR; $c++) {
        $newVals[$i, $c-1] = $oldVals[$srcIdx+1, $c]
    }
    for ($c = 1; $c -le $colCountSY; $c++) {
        $newFormulas[$i, $c-1] = $oldFormulas[$srcIdx+1, $c]
    }
}

# ---------------------------------------------------------------------
# 2. Write the reordered blocks back.
# ---------------------------------------------------------------------
$ws.Range("A5:R33").Value2 = $newVals
$ws.Range("S5:Y33").Formula = $newFormulas

# ---------------------------------------------------------------------
# 3. Bump the "Förändrad" column (C) to 46079 for every data row.
# ---------------------------------------------------------------------
$ws.Range("C2:C33").Value = 46079
